$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnas"
$ws.Range("C2").Value = "Vipr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 84.76851766666665
$ws.Range("H2").Value = 254.305553
$ws.Range("I2").Value = 0.2571740874301185
$ws.Range("J2").Value = 0.2571740874301185
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 1.802565333333333
$ws.Range("N2").Value = 5.407696
$ws.Range("O2").Value = 0.3083214936279398
$ws.Range("P2").Value = 0.3083214936279398
$ws.Range("Q2").Value = 152.8007913039875
$ws.Range("R2").Value = 1375.207121735888
$ws.Range("S2").Value = 0.0792922987588565
$ws.Range("T2").Value = 0.07929229875885653

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnas"
$ws.Range("C3").Value = "Vipr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 84.76851766666665
$ws.Range("H3").Value = 254.305553
$ws.Range("I3").Value = 0.2571740874301185
$ws.Range("J3").Value = 0.2571740874301185
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.06345833333333332
$ws.Range("N3").Value = 0.190375
$ws.Range("O3").Value = 0.01085429069041955
$ws.Range("P3").Value = 0.01085429069041955
$ws.Range("Q3").Value = 5.379268850263887
$ws.Range("R3").Value = 48.41341965237499
$ws.Range("S3").Value = 0.002791442303009878
$ws.Range("T3").Value = 0.002791442303009879

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Gnas"
$ws.Range("C4").Value = "Vipr1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 84.76851766666665
$ws.Range("H4").Value = 254.305553
$ws.Range("I4").Value = 0.2571740874301185
$ws.Range("J4").Value = 0.2571740874301185
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.980358666666667
$ws.Range("N4").Value = 11.941076
$ws.Range("O4").Value = 0.6808242156816406
$ws.Range("P4").Value = 0.6808242156816406
$ws.Range("Q4").Value = 337.4091039550031
$ws.Range("R4").Value = 3036.681935595028
$ws.Range("S4").Value = 0.1750903463682521
$ws.Range("T4").Value = 0.1750903463682521

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gnas"
$ws.Range("C5").Value = "Vipr1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 133.6830266666667
$ws.Range("H5").Value = 401.04908
$ws.Range("I5").Value = 0.4055728628296552
$ws.Range("J5").Value = 0.4055728628296552
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 1.802565333333333
$ws.Range("N5").Value = 5.407696
$ws.Range("O5").Value = 0.3083214936279398
$ws.Range("P5").Value = 0.3083214936279398
$ws.Range("Q5").Value = 240.9723895244089
$ws.Range("R5").Value = 2168.75150571968
$ws.Range("S5").Value = 0.1250468308425988
$ws.Range("T5").Value = 0.1250468308425988

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Gnas"
$ws.Range("C6").Value = "Vipr1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 133.6830266666667
$ws.Range("H6").Value = 401.04908
$ws.Range("I6").Value = 0.4055728628296552
$ws.Range("J6").Value = 0.4055728628296552
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.06345833333333332
$ws.Range("N6").Value = 0.190375
$ws.Range("O6").Value = 0.01085429069041955
$ws.Range("P6").Value = 0.01085429069041955
$ws.Range("Q6").Value = 8.483302067222223
$ws.Range("R6").Value = 76.34971860499999
$ws.Range("S6").Value = 0.004402205749298731
$ws.Range("T6").Value = 0.004402205749298732

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Gnas"
$ws.Range("C7").Value = "Vipr1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 133.6830266666667
$ws.Range("H7").Value = 401.04908
$ws.Range("I7").Value = 0.4055728628296552
$ws.Range("J7").Value = 0.4055728628296552
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.980358666666667
$ws.Range("N7").Value = 11.941076
$ws.Range("O7").Value = 0.6808242156816406
$ws.Range("P7").Value = 0.6808242156816406
$ws.Range("Q7").Value = 532.1063937788979
$ws.Range("R7").Value = 4788.957544010081
$ws.Range("S7").Value = 0.2761238262377576
$ws.Range("T7").Value = 0.2761238262377576

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Gnas"
$ws.Range("C8").Value = "Vipr1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 111.16377
$ws.Range("H8").Value = 333.49131
$ws.Range("I8").Value = 0.3372530497402263
$ws.Range("J8").Value = 0.3372530497402264
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 1.802565333333333
$ws.Range("N8").Value = 5.407696
$ws.Range("O8").Value = 0.3083214936279398
$ws.Range("P8").Value = 0.3083214936279398
$ws.Range("Q8").Value = 200.37995812464
$ws.Range("R8").Value = 1803.41962312176
$ws.Range("S8").Value = 0.1039823640264844
$ws.Range("T8").Value = 0.1039823640264845

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Gnas"
$ws.Range("C9").Value = "Vipr1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 111.16377
$ws.Range("H9").Value = 333.49131
$ws.Range("I9").Value = 0.3372530497402263
$ws.Range("J9").Value = 0.3372530497402264
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.06345833333333332
$ws.Range("N9").Value = 0.190375
$ws.Range("O9").Value = 0.01085429069041955
$ws.Range("P9").Value = 0.01085429069041955
$ws.Range("Q9").Value = 7.054267571249999
$ws.Range("R9").Value = 63.48840814125
$ws.Range("S9").Value = 0.003660642638110939
$ws.Range("T9").Value = 0.00366064263811094

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Gnas"
$ws.Range("C10").Value = "Vipr1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 111.16377
$ws.Range("H10").Value = 333.49131
$ws.Range("I10").Value = 0.3372530497402263
$ws.Range("J10").Value = 0.3372530497402264
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.980358666666667
$ws.Range("N10").Value = 11.941076
$ws.Range("O10").Value = 0.6808242156816406
$ws.Range("P10").Value = 0.6808242156816406
$ws.Range("Q10").Value = 442.47167533884
$ws.Range("R10").Value = 3982.24507804956
$ws.Range("S10").Value = 0.2296100430756309
$ws.Range("T10").Value = 0.2296100430756309
